# Update Sheet1 "Goals by each team" / "Annual Goals Count" boolean-style
# columns (B, D, F, G for rows 2-5) from numeric 1/0 flags to textual
# "yes"/"no" flags.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (was: B2=1, D2=0, F2=0, G2=0)
$ws.Range("B2").Value = "yes"
$ws.Range("D2").Value = "no"
$ws.Range("F2").Value = "no"
$ws.Range("G2").Value = "no"

# Row 3 (was: B3=0, D3=1, F3=0, G3=0)
$ws.Range("B3").Value = "no"
$ws.Range("D3").Value = "yes"
$ws.Range("F3").Value = "no"
$ws.Range("G3").Value = "no"

# Row 4 (was: B4=0, D4=0, F4=0, G4=0)
$ws.Range("B4").Value = "no"
$ws.Range("D4").Value = "no"
$ws.Range("F4").Value = "no"
$ws.Range("G4").Value = "no"

# Row 5 (was: B5=0, D5=0, F5=0, G5=0)
$ws.Range("B5").Value = "no"
$ws.Range("D5").Value = "no"
$ws.Range("F5").Value = "no"
$ws.Range("G5").Value = "no"

# The author's last selection ended up on E7 (empty cell below the table)
# instead of the previous G5.
[void]$ws.Range("E7").Select()
